# Prototype for NiCd allocated process
# This script edits the "batt_lci" worksheet:
#  - Updates the "code" (UUID) values for several activities
#  - Adds a new "source" column (F) to the waste-treatment activity's exchange table
#  - Inserts a brand-new exchange row describing a "market for copper smelting
#    facility" input, and updates the (now shifted) electricity-market exchange
#    row's amount and adds its own "source" value

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update activity "code" (UUID) values that changed ---------------
$ws.Range("B4").Value  = "70ee3a39-0833-4d35-b404-a39e9bff7b97"
$ws.Range("B14").Value = "098fcc52-e3ec-490a-beb3-c05c7d07f218"
$ws.Range("B24").Value = "4953ef88-d5d5-4523-9b2e-b7e66cc3ae2e"

# --- 2. Add the new "source" column header to the waste activity table ---
$ws.Range("F30").Value = "source"

# --- 3. Insert a new row right before the electricity-market exchange row,
#        which pushes everything below it down by one row ---------------
$ws.Rows.Item(32).Insert()

# --- 4. Populate the brand new row 32 (copper smelting facility input) ---
$ws.Range("A32").Value = "market for copper smelting facility"
$ws.Range("B32").Value = 0.0000000005
$ws.Range("C32").Value = "GLO"
$ws.Range("D32").Value = "unit"
$ws.Range("E32").Value = "technosphere"
$ws.Range("F32").Value = "ecoinvent treatment of Ni-metal hybride battery"

# --- 5. Update the electricity-market row (now shifted to row 33) --------
$ws.Range("B33").Value = 0.31
$ws.Range("F33").Value = "ecoinvent treatment of Ni-metal hybride battery"

# --- 6. Update the "code" (UUID) values for the two remaining activities -
$ws.Range("B40").Value = "6ab66874-3b7f-4c77-8601-719cb6ea6438"
$ws.Range("B50").Value = "053574a8-569f-4fdd-b61a-4053f9a313f5"
